# Remove the example row from the "Costs" sheet of the IHR Costing Tool
# Line Item Export template.
#
# The example data that lived in row 2 (columns A:Q) is cleared out, which
# leaves the formulas in U2/V2 evaluating against blank inputs. The row's
# height collapses back to the sheet default (no more wrapped example
# text), and the view/selection that had been left scrolled over to the
# example values is reset to the left-hand side of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Costs")

# Clear the example values out of row 2 (A2:Q2). R2:T2 are already empty.
$ws.Range("A2:Q2").ClearContents()

# The row no longer needs the taller, wrapped-text height used for the
# long example descriptions; it settles back down to (near) the sheet's
# default single-line height once the wrapped text is gone.
$ws.Rows.Item(2).RowHeight = 15.75

# Reset the view back to the left of the sheet instead of being scrolled
# over to where the example values used to be.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 10
$ws.Range("S1").Select()
